# Fill out the rubric checkboxes (column D) with "x" for the criteria that
# have been completed so far, and fix the selection to reflect where the
# editor left off (E26) instead of the old L10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(7,8,9,10,11,13,14,15,16,17,18,19,20,21,22,24,26,27,28,29,30,31,32,33,34,35,36,38,39,41,42,49,50,51,52,56,57,58,59,60,61,62,63,64)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = "x"
}

[void]$ws.Range("E26").Select()
